$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")

# New Policia data for the 2025 sheet (Hombres / Mujeres counts per age group)
$data = @(
    @(4, 25),
    @(9, 99),
    @(3, 41),
    @(5, 50),
    @(12, 57),
    @(26, 55),
    @(27, 85),
    @(44, 88),
    @(14, 58),
    @(17, 57),
    @(7, 40),
    @(7, 37),
    @(4, 30),
    @(2, 73),
    @(67, 371)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

# Make the 2025 sheet active with its own selection, matching the saved view state
$ws.Activate()
$ws.Range("E9").Select()
